$d = $word.ActiveDocument

# Docx writer: Use different style for block quotes in notes.
# Add a new "Footnote Block Text" paragraph style, based on "Footnote Text"
# (so footnote-note block quotes can get their own font size), mirroring the
# existing "Block Text" style's spacing/indent but chaining back into
# "Footnote Text" for BaseStyle/Next.
$newStyle = $d.Styles.Add("Footnote Block Text", 1)
$newStyle.BaseStyle = "Footnote Text"
$newStyle.NextParagraphStyle = "Footnote Text"
$newStyle.Priority = 9
$newStyle.UnhideWhenUsed = $true
$newStyle.QuickStyle = $true

# <w:pPr><w:spacing w:after="100" w:before="100"/><w:ind w:firstLine="0" w:left="480" w:right="480"/></w:pPr>
$newStyle.ParagraphFormat.SpaceBefore = 5
$newStyle.ParagraphFormat.SpaceAfter = 5
$newStyle.ParagraphFormat.FirstLineIndent = 0
$newStyle.ParagraphFormat.LeftIndent = 24
$newStyle.ParagraphFormat.RightIndent = 24
